$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "productivity" (B) and "ratio" (D) columns previously carried a custom
# 2-decimal-place number format. Drop that formatting so these columns share
# the same plain style as the "year" (A) and "empl_cost" (C) columns, by
# copying the formats already used in those neighboring columns.
# -4122 == xlPasteFormats
$xlPasteFormats = -4122
$ws.Range("A2").Copy()
$ws.Range("B2:B7").PasteSpecial($xlPasteFormats)
$ws.Range("C2").Copy()
$ws.Range("D2:D7").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Append the new 2021 data row to the table.
$ws.Range("A8").Value = 2021
$ws.Range("B8").Value = 0.6
$ws.Range("C8").Value = 126
$ws.Range("D8").Value = 3.2

# Match the active cell selection recorded in the saved workbook view.
$ws.Range("H12").Select()
